$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 585.7293526225917
$ws.Range("D2").Value = 4185.871936422012
